$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 18:52"

# Re-rank rows 21-31: Ecuador jumps above Peru..Mexico which shift down one place
$ws.Range("A21").Value = "Ecuador"
$ws.Range("A22").Value = "Peru"
$ws.Range("A23").Value = "Irlanda"
$ws.Range("A24").Value = "Suecia"
$ws.Range("A25").Value = "Arabia Saudita"
$ws.Range("A26").Value = "Austria"
$ws.Range("A27").Value = "Israel"
$ws.Range("A28").Value = "Japon"
$ws.Range("A29").Value = "Chile"
$ws.Range("A30").Value = "Singapur"
$ws.Range("A31").Value = "Mexico"

# Update numeric stats per the source refresh
# Row 4
$ws.Range("B4").Value = 894475
$ws.Range("C4").Value = 8033
$ws.Range("D4").Value = 90201
$ws.Range("E4").Value = 753335
$ws.Range("F4").Value = 15042
$ws.Range("G4").Value = 703
$ws.Range("H4").Value = 50939

# Row 6
$ws.Range("B6").Value = 192994
$ws.Range("C6").Value = 3021
$ws.Range("D6").Value = 60498
$ws.Range("E6").Value = 106527
$ws.Range("F6").Value = 2173
$ws.Range("G6").Value = 420
$ws.Range("H6").Value = 25969

# Row 8
$ws.Range("B8").Value = 154111
$ws.Range("C8").Value = 982
$ws.Range("D8").Value = 106800
$ws.Range("E8").Value = 41679
$ws.Range("F8").Value = 2908
$ws.Range("G8").Value = 57
$ws.Range("H8").Value = 5632

# Row 10
$ws.Range("B10").Value = 104912
$ws.Range("C10").Value = 3122
$ws.Range("D10").Value = 21737
$ws.Range("E10").Value = 80575
$ws.Range("F10").Value = 1790
$ws.Range("G10").Value = 109
$ws.Range("H10").Value = 2600

# Row 14
$ws.Range("B14").Value = 51073
$ws.Range("C14").Value = 1581
$ws.Range("D14").Value = 26573
$ws.Range("E14").Value = 21093
$ws.Range("F14").Value = 8318
$ws.Range("G14").Value = 94
$ws.Range("H14").Value = 3407

# Row 21
$ws.Range("B21").Value = 22719
$ws.Range("C21").Value = 11536
$ws.Range("D21").Value = 1366
$ws.Range("E21").Value = 20777
$ws.Range("F21").Value = 127
$ws.Range("G21").Value = 16
$ws.Range("H21").Value = 576

# Row 22
$ws.Range("B22").Value = 20914
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 7422
$ws.Range("E22").Value = 12920
$ws.Range("F22").Value = 396
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 572

# Row 23
$ws.Range("B23").Value = 17607
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 9233
$ws.Range("E23").Value = 7580
$ws.Range("F23").Value = 147
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 794

# Row 24
$ws.Range("B24").Value = 17567
$ws.Range("C24").Value = 812
$ws.Range("D24").Value = 550
$ws.Range("E24").Value = 14865
$ws.Range("F24").Value = 547
$ws.Range("G24").Value = 131
$ws.Range("H24").Value = 2152

# Row 25
$ws.Range("B25").Value = 15102
$ws.Range("C25").Value = 1172
$ws.Range("D25").Value = 2049
$ws.Range("E25").Value = 12926
$ws.Range("F25").Value = 93
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 127

# Row 26
$ws.Range("B26").Value = 15071
$ws.Range("C26").Value = 69
$ws.Range("D26").Value = 11872
$ws.Range("E26").Value = 2669
$ws.Range("F26").Value = 156
$ws.Range("G26").Value = 8
$ws.Range("H26").Value = 530

# Row 27
$ws.Range("B27").Value = 14882
$ws.Range("C27").Value = 79
$ws.Range("D27").Value = 5685
$ws.Range("E27").Value = 9004
$ws.Range("F27").Value = 139
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 193

# Row 28
$ws.Range("B28").Value = 12368
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 1494
$ws.Range("E28").Value = 10546
$ws.Range("F28").Value = 259
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 328

# Row 29
$ws.Range("B29").Value = 12306
$ws.Range("C29").Value = 494
$ws.Range("D29").Value = 6327
$ws.Range("E29").Value = 5805
$ws.Range("F29").Value = 408
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 174

# Row 30
$ws.Range("B30").Value = 12075
$ws.Range("C30").Value = 897
$ws.Range("D30").Value = 924
$ws.Range("E30").Value = 11139
$ws.Range("F30").Value = 26
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 12

# Row 31
$ws.Range("B31").Value = 11633
$ws.Range("C31").Value = 1089
$ws.Range("D31").Value = 2627
$ws.Range("E31").Value = 7937
$ws.Range("F31").Value = 378
$ws.Range("G31").Value = 99
$ws.Range("H31").Value = 1069

# Row 135
$ws.Range("B135").Value = 144
$ws.Range("C135").Value = 12
$ws.Range("D135").Value = 9
$ws.Range("E135").Value = 130
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 5
